# "Added feature to copy and move cards"
#
# The sheet is a flattened Trello-style board export: column A holds each
# card/list's unique id, column B holds its display name. Two extra "Week 1"
# cards were added (the board now shows 4 cards copied/moved into the
# "Week 1" list instead of 1), which re-generated every id below it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make room for the two new "Week 1" cards - shift rows 3..18 down to 5..20.
$ws.Rows.Item(3).Insert()
$ws.Rows.Item(3).Insert()

# Also widen the default column width to fit the new card text.
$ws.StandardWidth = 31.432656

$rows = @(
    @{ Row=2;  A="62b6857248c19f7c8b46ebf3"; B="Week 1" },
    @{ Row=3;  A="62b68570affd15341d40cf22"; B="Week 1" },
    @{ Row=4;  A="62b6856e9bc5355be91c56bd"; B="Week 1" },
    @{ Row=5;  A="62b6856cc0052979cfd086dd"; B="Week 1" },
    @{ Row=6;  A="62b67a3b2be36f4542da6bcb"; B="Week 1" },
    @{ Row=7;  A="62b67a3baaaf838f4daaceb0"; B="Week 2" },
    @{ Row=8;  A="62b67a3a2203915bd0454b2f"; B="Week 3" },
    @{ Row=9;  A="62b67a3a4512494794b52d86"; B="Week 4" },
    @{ Row=10; A="62b67a3a9d5f007785c4056f"; B="Week 5" },
    @{ Row=11; A="62b67a399d4dca039dc1a56f"; B="Week 6" },
    @{ Row=12; A="62b67a3981ba16144fecfc34"; B="Week 7" },
    @{ Row=13; A="62b67a395e18168c34b79d8b"; B="Week 8" },
    @{ Row=14; A="62b67a3814716f4f28fe2431"; B="Week 9" },
    @{ Row=15; A="62b67a3810ce2a19b86075d4"; B="Week 10" },
    @{ Row=16; A="62b67a38e7c29b8ec3dc7d1b"; B="Week 11" },
    @{ Row=17; A="62b67a3821622e7ae813dbef"; B="Week 12" },
    @{ Row=18; A="62b67a354ca480872e59caea"; B="To Do" },
    @{ Row=19; A="62b67a354ca480872e59caeb"; B="Doing" },
    @{ Row=20; A="62b67a354ca480872e59caec"; B="Done" }
)

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 1).Value = $r.A
    $ws.Cells.Item($r.Row, 2).Value = $r.B
}

Write-Output "applied card copy/move edit"
